$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E10:E40").Copy()
$ws.Range("F10:H40").PasteSpecial(-4122)
